$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 15386368
$ws.Range("I100").Value = 1086.3334
$ws.Range("K100").Value = 1086.3334
$ws.Range("M100").Value = -545.3334

# Row 121
$ws.Range("H121").Value = 1797.5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1797.5
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").Value = 5392.5
$ws.Range("N121").Value = -8886.5

# Row 132
$ws.Range("H132").Value = 2021717.5
$ws.Range("I132").Value = 1103.4791
$ws.Range("J132").Value = 15877357
$ws.Range("K132").Value = 3310.4373
$ws.Range("L132").Value = 47632071
$ws.Range("M132").Value = -780.4373000000001
$ws.Range("N132").Value = -47637131

# Row 137
$ws.Range("H137").Value = 11584064
$ws.Range("I137").Value = 908.87177
$ws.Range("J137").Value = 41700264
$ws.Range("K137").Value = 2726.61531
$ws.Range("L137").Value = 125100792
$ws.Range("M137").Value = -176.6153100000001
$ws.Range("N137").Value = -125105892

# Row 138
$ws.Range("H138").Value = 2281.4583
$ws.Range("I138").Value = 1475.1887
$ws.Range("J138").Value = 4530.5264
$ws.Range("K138").Value = 4425.5661
$ws.Range("L138").Value = 13591.5792
$ws.Range("M138").Value = 714.4339
$ws.Range("N138").Value = -23871.5792

# Row 141
$ws.Range("H141").Value = 1361.7954
$ws.Range("I141").Value = 939.8421
$ws.Range("J141").Value = 4034.1667
$ws.Range("K141").Value = 2819.5263
$ws.Range("L141").Value = 12102.5001
$ws.Range("M141").Value = 2360.4737
$ws.Range("N141").Value = -22462.5001


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 43810970
$ws.Range("I74").Value = 38462212
$ws.Range("J74").Value = 59262920
$ws.Range("K74").Value = 38462212
$ws.Range("L74").Value = 59262920
$ws.Range("M74").Value = -38461338
$ws.Range("N74").Value = -59264668

# Row 77
$ws.Range("H77").Value = 43810970
$ws.Range("I77").Value = 38462212
$ws.Range("J77").Value = 59262920
$ws.Range("K77").Value = 192311060
$ws.Range("L77").Value = 296314600
$ws.Range("M77").Value = -192306692
$ws.Range("N77").Value = -296323336

# Row 88
$ws.Range("H88").Value = 5769.231
$ws.Range("I88").Value = 2300
$ws.Range("J88").Value = 6810
$ws.Range("K88").Value = 2300
$ws.Range("L88").Value = 6810
$ws.Range("M88").Value = -1894
$ws.Range("N88").Value = -7622

# Row 91
$ws.Range("H91").Value = 5769.231
$ws.Range("I91").Value = 2300
$ws.Range("J91").Value = 6810
$ws.Range("K91").Value = 2300
$ws.Range("L91").Value = 6810
$ws.Range("M91").Value = -896
$ws.Range("N91").Value = -9618

# Row 102
$ws.Range("H102").Value = 2179.4
$ws.Range("I102").Value = 2283.3333
$ws.Range("J102").Value = 1763.6666
$ws.Range("K102").Value = 2283.3333
$ws.Range("L102").Value = 1763.6666
$ws.Range("M102").Value = -661.3332999999998
$ws.Range("N102").Value = -5007.6666


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1872.67
$ws.Range("I86").Value = 1891.4124
$ws.Range("J86").Value = 1266.6666
$ws.Range("K86").Value = 1891.4124
$ws.Range("L86").Value = 1266.6666
$ws.Range("M86").Value = -768.4123999999999
$ws.Range("N86").Value = -3512.6666

# Row 89
$ws.Range("H89").Value = 1872.67
$ws.Range("I89").Value = 1891.4124
$ws.Range("J89").Value = 1266.6666
$ws.Range("K89").Value = 9457.062
$ws.Range("L89").Value = 6333.333000000001
$ws.Range("M89").Value = -3841.062
$ws.Range("N89").Value = -17565.333

# Row 94
$ws.Range("H94").Value = 1185.25
$ws.Range("I94").Value = 855.5217
$ws.Range("J94").Value = 2702
$ws.Range("K94").Value = 855.5217
$ws.Range("L94").Value = 2702
$ws.Range("M94").Value = -404.5217
$ws.Range("N94").Value = -3604


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 2287.1428
$ws.Range("I14").Value = 1005
$ws.Range("J14").Value = 2800
$ws.Range("K14").Value = 1005
$ws.Range("L14").Value = 2800
$ws.Range("M14").Value = -835
$ws.Range("N14").Value = -3140

# Row 21
$ws.Range("H21").Value = 4013
$ws.Range("I21").Value = 4013
$ws.Range("K21").Value = 4013
$ws.Range("M21").Value = -3778


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 855.04346
$ws.Range("I131").Value = 348.17648
$ws.Range("J131").Value = 1152.1724
$ws.Range("K131").Value = 1044.52944
$ws.Range("L131").Value = 3456.5172
$ws.Range("M131").Value = 3995.47056
$ws.Range("N131").Value = -13536.5172


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 34817.55
$ws.Range("I16").Value = 38719.77
$ws.Range("J16").Value = 998.3333
$ws.Range("K16").Value = 38719.77
$ws.Range("L16").Value = 998.3333
$ws.Range("M16").Value = -38549.77
$ws.Range("N16").Value = -1338.3333

# Row 68
$ws.Range("H68").Value = 2659.85
$ws.Range("I68").Value = 2406.6667
$ws.Range("J68").Value = 3419.4
$ws.Range("K68").Value = 2406.6667
$ws.Range("L68").Value = 3419.4
$ws.Range("M68").Value = -1657.6667
$ws.Range("N68").Value = -4917.4

# Row 71
$ws.Range("H71").Value = 2659.85
$ws.Range("I71").Value = 2406.6667
$ws.Range("J71").Value = 3419.4
$ws.Range("K71").Value = 12033.3335
$ws.Range("L71").Value = 17097
$ws.Range("M71").Value = -8289.333500000001
$ws.Range("N71").Value = -24585

# Row 136
$ws.Range("H136").Value = 9261841
$ws.Range("I136").Value = 12347788
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 37043364
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -37040814
$ws.Range("N136").Value = -17100


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 31267388
$ws.Range("I62").Value = 83369800
$ws.Range("J62").Value = 5939.9
$ws.Range("K62").Value = 83369800
$ws.Range("L62").Value = 5939.9
$ws.Range("M62").Value = -83369176
$ws.Range("N62").Value = -7187.9

# Row 65
$ws.Range("H65").Value = 31267388
$ws.Range("I65").Value = 83369800
$ws.Range("J65").Value = 5939.9
$ws.Range("K65").Value = 416849000
$ws.Range("L65").Value = 29699.5
$ws.Range("M65").Value = -416845880
$ws.Range("N65").Value = -35939.5

# Row 136
$ws.Range("H136").Value = 2602.5874
$ws.Range("I136").Value = 633.5208
$ws.Range("J136").Value = 8903.6
$ws.Range("K136").Value = 1900.5624
$ws.Range("L136").Value = 26710.8
$ws.Range("M136").Value = 649.4376
$ws.Range("N136").Value = -31810.8

